# Refresh the cryptos list cells to the latest scraped values/prices.
# (GitHub Actions scheduled update, mirrors the upstream commit diff.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.374.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.82%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.771.12"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.49%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.43%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.767.86"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.56%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("E9").Value = "  -1.39%  "

# Row 10
$ws.Range("E10").Value = "  -2.11%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.60"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.12%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.486"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.01"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.68%  "

# Row 14
$ws.Range("E14").Value = "  -2.64%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.395.56"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.767.89"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.436.48"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.76%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.57"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.38%  "

# Row 19
$ws.Range("E19").Value = "  -3.30%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.89"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.03%  "

# Row 21
$ws.Range("E21").Value = "  -3.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.53%  "

# Row 23
$ws.Range("E23").Value = "  +0.43%  "

# Row 24
$ws.Range("E24").Value = "  -0.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.84%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.85"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.91%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000136"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.94%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.84%  "

# Row 29
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.53"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.50%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.99"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.76"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.84%  "

# Row 34
$ws.Range("E34").Value = "  -0.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.09%  "

# Row 36
$ws.Range("E36").Value = "  -2.08%  "

# Row 37
$ws.Range("E37").Value = "  -1.30%  "

# Row 38
$ws.Range("E38").Value = "  +4.64%  "

# Row 39
$ws.Range("E39").Value = "  +2.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "455.64"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.08"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.20%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.70%  "

# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.78"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.33%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.67"
$ws.Range("D44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.58"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.67%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.956.79"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.72%  "

# Row 47
$ws.Range("E47").Value = "  -0.38%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.11"
$ws.Range("D50").ClearFormats()

# Row 51
$ws.Range("E51").Value = "  +0.30%  "
